# Insert a new data row at row 298 (pushing the existing rows 298-328 down
# to 299-329) and populate it with a new Espinaca price record for
# Terminal La Palmera de La Serena, matching the rest of the dataset's
# static columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 298:328 down one row, leaving row 298 empty for the new record.
$ws.Rows.Item(298).Insert()

$ws.Range("A298").Value = 8
$ws.Range("B298").Value = "Terminal La Palmera de La Serena"
$ws.Range("C298").Value = "Coquimbo"
$ws.Range("D298").Value = 44858
$ws.Range("E298").Value = 4
$ws.Range("F298").Value = 100112012
$ws.Range("G298").Value = "Espinaca"
$ws.Range("H298").Value = "Sin especificar"
$ws.Range("I298").Value = "Primera"
$ws.Range("J298").Value = 2400
$ws.Range("K298").Value = 450
$ws.Range("L298").Value = 500
$ws.Range("M298").Value = 475
$ws.Range("N298").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O298").Value = "Provincia del Elquí"
$ws.Range("P298").Value = 950
$ws.Range("Q298").Value = 0.5
$ws.Range("R298").Value = "Hortaliza"
